$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '72.959.92'
$ws.Range('E2').Value = '  +2.12%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.996.18'
$ws.Range('E3').Value = '  +0.54%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '623.03'
$ws.Range('E5').Value = '  +15.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '163.50'
$ws.Range('E6').Value = '  +8.95%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.688'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  +1.74%  '
$ws.Range('E10').Value = '  +0.54%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.46'
$ws.Range('E11').Value = '  -1.64%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000320'
$ws.Range('E12').Value = '  -0.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.10'
$ws.Range('E13').Value = '  +3.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.627.92'
$ws.Range('E14').Value = '  +0.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.997.16'
$ws.Range('E15').Value = '  +0.61%  '
$ws.Range('E16').Value = '  +7.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.20'
$ws.Range('E17').Value = '  +0.95%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.72'
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.635.83'
$ws.Range('E20').Value = '  +1.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '441.40'
$ws.Range('E21').Value = '  +2.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.99'
$ws.Range('E22').Value = '  +17.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '97.05'
$ws.Range('E23').Value = '  -0.44%  '
$ws.Range('E24').Value = '  -3.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '14.49'
$ws.Range('E25').Value = '  -1.25%  '
$ws.Range('E26').Value = '  +4.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.35'
$ws.Range('E27').Value = '  -0.92%  '
$ws.Range('B28').Value = 'Filecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.63'
$ws.Range('E28').Value = '  -2.22%  '
$ws.Range('B29').Value = 'LEO'
$ws.Range('C29').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.98'
$ws.Range('E29').Value = '  +1.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.48'
$ws.Range('E30').Value = '  -1.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.74'
$ws.Range('E31').Value = '  -2.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '14.00'
$ws.Range('E32').Value = '  +3.90%  '
$ws.Range('E33').Value = '  -0.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '72.67'
$ws.Range('E34').Value = '  +10.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '48.21'
$ws.Range('E35').Value = '  -6.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '639.63'
$ws.Range('E36').Value = '  -6.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0₃0902'
$ws.Range('E37').Value = '  +9.06%  '
$ws.Range('E38').Value = '  -0.89%  '
$ws.Range('E39').Value = '  -0.89%  '
$ws.Range('E40').Value = '  -0.85%  '
$ws.Range('E41').Value = '  -0.22%  '
$ws.Range('E42').Value = '  +4.02%  '
$ws.Range('E43').Value = '  +0.12%  '
$ws.Range('E44').Value = '  +1.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.66'
$ws.Range('E45').Value = '  +2.66%  '
$ws.Range('E46').Value = '  +0.89%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.66'
$ws.Range('E47').Value = '  -0.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.43'
$ws.Range('E48').Value = '  +1.84%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.925.74'
$ws.Range('E49').Value = '  +11.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.09'
$ws.Range('E50').Value = '  +1.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.43'
